$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price"/"Volume(1h)" columns (D:E) store plain-looking numeric text
# (e.g. "205.02", "79.368.41", "  +3.62%  ") as literal strings, not numbers.
# Excel's COM layer auto-coerces a bare numeric-looking string assigned via
# .Value into a real Number, which would silently change the cell type and
# drop formatting like trailing zeros. Forcing the whole data range to Text
# format first keeps every assignment below as a literal string; the
# trailing ClearFormats() restores the original (unstyled) look, matching
# the workbook's existing formatting for these cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.Value = $value
}

Set-TextValue "D2" "79.368.41"
Set-TextValue "E2" "  +3.62%  "
Set-TextValue "D3" "3.187.62"
Set-TextValue "E3" "  +4.67%  "
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "205.02"
Set-TextValue "E5" "  +1.52%  "
Set-TextValue "D6" "633.77"
Set-TextValue "E6" "  +0.57%  "
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "E8" "  +10.63%  "
Set-TextValue "E9" "  +5.57%  "
Set-TextValue "D10" "3.186.75"
Set-TextValue "E10" "  +4.70%  "
Set-TextValue "D11" "0.586"
Set-TextValue "E11" "  +33.29%  "
Set-TextValue "E12" "  +3.00%  "
Set-TextValue "D13" "5.53"
Set-TextValue "E13" "  +7.48%  "
Set-TextValue "D14" "3.774.01"
Set-TextValue "E14" "  +4.82%  "
Set-TextValue "D15" "0.0000226"
Set-TextValue "E15" "  +16.22%  "
Set-TextValue "D16" "31.75"
Set-TextValue "E16" "  +7.31%  "
Set-TextValue "D17" "79.322.61"
Set-TextValue "E17" "  +3.67%  "
Set-TextValue "D18" "3.196.95"
Set-TextValue "E18" "  +5.23%  "
Set-TextValue "D19" "14.48"
Set-TextValue "E19" "  +7.35%  "
Set-TextValue "D20" "3.07"
Set-TextValue "E20" "  +33.49%  "
Set-TextValue "D21" "9.09"
Set-TextValue "E21" "  +0.12%  "
Set-TextValue "D22" "430.67"
Set-TextValue "E22" "  +14.30%  "
Set-TextValue "D23" "5.01"
Set-TextValue "E23" "  +14.46%  "
Set-TextValue "B24" "WrappedeETH"
Set-TextValue "C24" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D24" "3.355.68"
Set-TextValue "E24" "  +5.39%  "
Set-TextValue "B25" "NEARProtocol"
Set-TextValue "C25" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D25" "4.78"
Set-TextValue "E25" "  +8.97%  "
Set-TextValue "B26" "Aptos"
Set-TextValue "C26" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D26" "11.18"
Set-TextValue "E26" "  +11.86%  "
Set-TextValue "D27" "76.78"
Set-TextValue "E27" "  +4.22%  "
Set-TextValue "D28" "0.997"
Set-TextValue "E28" "  -0.37%  "
Set-TextValue "D29" "0.0000118"
Set-TextValue "E29" "  +5.45%  "
Set-TextValue "B30" "Binance-PegBSC-USD"
Set-TextValue "C30" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D30" "0.998"
Set-TextValue "E30" "  -0.08%  "
Set-TextValue "B31" "InternetComputer(DFINITY)"
Set-TextValue "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "8.98"
Set-TextValue "E31" "  +7.49%  "
Set-TextValue "D32" "1.48"
Set-TextValue "E32" "  +4.04%  "
Set-TextValue "D33" "523.42"
Set-TextValue "E33" "  +1.68%  "
Set-TextValue "D34" "2.00"
Set-TextValue "E34" "  +2.17%  "
Set-TextValue "D35" "0.141"
Set-TextValue "E35" "  +25.96%  "
Set-TextValue "D36" "22.89"
Set-TextValue "E36" "  +9.19%  "
Set-TextValue "E37" "  +10.43%  "
Set-TextValue "D38" "0.999"
Set-TextValue "E38" "  +0.04%  "
Set-TextValue "D39" "0.402"
Set-TextValue "E39" "  +4.08%  "
Set-TextValue "D40" "165.24"
Set-TextValue "E40" "  +1.09%  "
Set-TextValue "B41" "WhiteBITCoin"
Set-TextValue "C41" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D41" "20.03"
Set-TextValue "E41" "  +0.05%  "
Set-TextValue "B42" "Aave"
Set-TextValue "C42" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D42" "193.19"
Set-TextValue "E42" "  +2.48%  "
Set-TextValue "D44" "5.47"
Set-TextValue "E44" "  +4.94%  "
Set-TextValue "E45" "  +11.25%  "
Set-TextValue "E46" "  +7.27%  "
Set-TextValue "D47" "1.32"
Set-TextValue "E47" "  +3.42%  "
Set-TextValue "D48" "42.96"
Set-TextValue "E48" "  +2.16%  "
Set-TextValue "D49" "2.54"
Set-TextValue "E49" "  +3.94%  "
Set-TextValue "D50" "25.62"
Set-TextValue "E50" "  +13.25%  "
Set-TextValue "D51" "0.633"
Set-TextValue "E51" "  +4.21%  "

$dataRange.ClearFormats()
